$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "29.662.05"
$ws.Range("E2").Value = "  +2.75%  "

Set-TextValue "D3" "1.862.96"
$ws.Range("E3").Value = "  +2.15%  "

Set-TextValue "D4" "0.9991"

Set-TextValue "D5" "245.86"
$ws.Range("E5").Value = "  +2.92%  "

Set-TextValue "D6" "0.7006"
$ws.Range("E6").Value = "  +1.63%  "

Set-TextValue "D7" "0.9998"

Set-TextValue "D8" "0.07749"
$ws.Range("E8").Value = "  +1.90%  "

Set-TextValue "D9" "0.3073"
$ws.Range("E9").Value = "  +2.05%  "

Set-TextValue "D10" "23.70"
$ws.Range("E10").Value = "  +1.50%  "

Set-TextValue "D11" "0.07785"
$ws.Range("E11").Value = "  +0.85%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D12" "5.167"
$ws.Range("E12").Value = "  +2.56%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.859.91"
$ws.Range("E13").Value = "  +1.80%  "

Set-TextValue "D14" "92.48"
$ws.Range("E14").Value = "  +2.87%  "

Set-TextValue "D15" "0.6942"
$ws.Range("E15").Value = "  +3.53%  "

Set-TextValue "D16" "6.588"
$ws.Range("E16").Value = "  +2.79%  "

Set-TextValue "D17" "29.646.13"
$ws.Range("E17").Value = "  +2.75%  "

Set-TextValue "D18" "0.000008372"
$ws.Range("E18").Value = "  +1.28%  "

Set-TextValue "D19" "2.113.75"
$ws.Range("E19").Value = "  +1.61%  "

Set-TextValue "D20" "242.50"
$ws.Range("E20").Value = "  -0.05%  "

Set-TextValue "D21" "12.78"
$ws.Range("E21").Value = "  +1.52%  "

Set-TextValue "D22" "1.000"
$ws.Range("E22").Value = "  +0.06%  "

$ws.Range("E23").Value = "  +3.68%  "

Set-TextValue "D24" "0.9999"
$ws.Range("E24").Value = "  -0.01%  "

Set-TextValue "D25" "0.1515"
$ws.Range("E25").Value = "  +3.00%  "

Set-TextValue "D26" "8.935"
$ws.Range("E26").Value = "  +2.76%  "

Set-TextValue "D27" "159.96"
$ws.Range("E27").Value = "  -0.17%  "

Set-TextValue "D28" "18.34"
$ws.Range("E28").Value = "  +1.19%  "

Set-TextValue "D29" "1.541"
$ws.Range("E29").Value = "  +0.99%  "

Set-TextValue "D30" "4.267"
$ws.Range("E30").Value = "  +2.10%  "

Set-TextValue "D31" "4.194"
$ws.Range("E31").Value = "  +1.63%  "

Set-TextValue "D32" "1.193"
$ws.Range("E32").Value = "  +0.31%  "

Set-TextValue "D33" "0.05113"
$ws.Range("E33").Value = "  +0.43%  "

Set-TextValue "D34" "0.7867"
$ws.Range("E34").Value = "  +4.91%  "

Set-TextValue "D35" "1.907"
$ws.Range("E35").Value = "  +5.45%  "

Set-TextValue "D36" "1.160"
$ws.Range("E36").Value = "  +1.85%  "

Set-TextValue "D37" "2.686"
$ws.Range("E37").Value = "  +0.09%  "

Set-TextValue "D38" "1.332.47"
$ws.Range("E38").Value = "  +11.03%  "

Set-TextValue "D39" "0.01879"
$ws.Range("E39").Value = "  +2.84%  "

Set-TextValue "D40" "2.736"
$ws.Range("E40").Value = "  +2.43%  "

Set-TextValue "D41" "0.9593"
$ws.Range("E41").Value = "  +4.98%  "

Set-TextValue "D42" "6.009"
$ws.Range("E42").Value = "  +14.99%  "

Set-TextValue "D43" "106.58"
$ws.Range("E43").Value = "  -1.27%  "

Set-TextValue "D44" "0.9995"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("E45").Value = "  +3.62%  "

Set-TextValue "D46" "9.774"
$ws.Range("E46").Value = "  +3.77%  "

Set-TextValue "D47" "2.012.63"
$ws.Range("E47").Value = "  +1.66%  "

Set-TextValue "D48" "0.5216"
$ws.Range("E48").Value = "  +1.27%  "

Set-TextValue "D49" "65.09"
$ws.Range("E49").Value = "  +4.77%  "

Set-TextValue "D51" "7.011"
$ws.Range("E51").Value = "  +2.30%  "
